# Generate Report for Handback
#
# This script updates the localization-status workbook so that the
# zh-cn and de-de sheets reflect a completed handback: the Status column
# changes from "Ready for handoff" to "Handed back: in sync with en-US",
# new "Latest Target File" (F) and "Latest Handback File" (G) hyperlink
# columns are populated, and the "Latest Handback DateTime" (H) values are
# filled in (no longer the zero-date placeholder).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# The Overview sheet mirrors the same "Status" text (shared string) used by
# the language sheets below. Update it too so the shared string is fully
# replaced everywhere rather than leaving a stray, now-unused copy behind.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

function Update-LangSheet($SheetName, $H2Value, $H3Value, $D3Value) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Capture the existing hyperlink addresses (keyed by cell reference) so
    # the new Target/Handback columns can reuse the same link targets as
    # the corresponding source (A/D) columns.
    $links = @{}
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        $links[$addr] = $h.Address
    }

    $addrA2 = $links['$A$2']
    $addrD2 = $links['$D$2']
    $addrA3 = $links['$A$3']
    $addrD3 = $links['$D$3']

    $dispA2 = $ws.Range("A2").Value2
    $dispD2 = $ws.Range("D2").Value2
    $dispA3 = $ws.Range("A3").Value2
    $dispD3 = $ws.Range("D3").Value2

    # Status column now reflects a finished handback.
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # New "Latest Target File" column (F) - same file / same link as the
    # "Source File Name" column (A).
    $ws.Hyperlinks.Add($ws.Range("F2"), $addrA2, "", "", $dispA2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $addrA3, "", "", $dispA3)

    # New "Latest Handback File" column (G) - same file / same link as the
    # "Latest Handoff File" column (D).
    $ws.Hyperlinks.Add($ws.Range("G2"), $addrD2, "", "", $dispD2)
    $ws.Hyperlinks.Add($ws.Range("G3"), $addrD3, "", "", $dispD3)

    # "Latest Handback DateTime" column (H) picks up real timestamps.
    $ws.Range("H2").Value = $H2Value
    $ws.Range("H3").Value = $H3Value

    if ($SheetName -eq "de-de") {
        # G3's displayed value diverges from its hyperlink's display text
        # in the target workbook (mirrors D3 below) - overwrite the cell
        # value after the hyperlink has been created so the stored
        # hyperlink display text is left untouched.
        $ws.Range("G3").Value = $D3Value
        $ws.Range("D3").Value = $D3Value
    }
}

Update-LangSheet "zh-cn" "2016-03-12 08:44:03" "2016-03-12 08:44:03" ""
Update-LangSheet "de-de" "4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.64c3bc10457d2ae7e5488accef9d1770522e41d0.de-de.xlf" "4e5fd6f6-47d7-4f37-b39a-5c1da5783b90.64c3bc10457d2ae7e5488accef9d1770522e41d0.de-de.xlf" "2016-03-12 08:44:09"

Write-Host "Handback report generated."
